$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.027123
$ws.Range("H2").Value = 0.081369
$ws.Range("I2").Value = 0.0960827240265261
$ws.Range("J2").Value = 0.09608272402652611
$ws.Range("M2").Value = 7.60267
$ws.Range("N2").Value = 22.80801
$ws.Range("O2").Value = 0.1369613544323865
$ws.Range("P2").Value = 0.1558153465164618
$ws.Range("Q2").Value = 0.20620721841
$ws.Range("R2").Value = 1.85586496569
$ws.Range("S2").Value = 0.01315962002022622
$ws.Range("T2").Value = 0.01497116293843874
$ws.Range("G3").Value = 0.027123
$ws.Range("H3").Value = 0.081369
$ws.Range("I3").Value = 0.0960827240265261
$ws.Range("J3").Value = 0.09608272402652611
$ws.Range("N3").Value = 81.75477699999999
$ws.Range("O3").Value = 0.4909347632361489
$ws.Range("P3").Value = 0.5585164557377457
$ws.Range("Q3").Value = 0.7391449388569998
$ws.Range("R3").Value = 6.652304449712999
$ws.Range("S3").Value = 0.04717034937104683
$ws.Range("T3").Value = 0.0536637824809233
$ws.Range("G4").Value = 0.027123
$ws.Range("H4").Value = 0.081369
$ws.Range("I4").Value = 0.0960827240265261
$ws.Range("J4").Value = 0.09608272402652611
$ws.Range("M4").Value = 0.3673663333333333
$ws.Range("N4").Value = 1.102099
$ws.Range("O4").Value = 0.00661806846623527
$ws.Range("P4").Value = 0.007529106554252038
$ws.Range("Q4").Value = 0.009964077058999998
$ws.Range("R4").Value = 0.089676693531
$ws.Range("S4").Value = 0.0006358820460299383
$ws.Range("T4").Value = 0.0007234170672185075
$ws.Range("G5").Value = 0.027123
$ws.Range("H5").Value = 0.081369
$ws.Range("I5").Value = 0.0960827240265261
$ws.Range("J5").Value = 0.09608272402652611
$ws.Range("M5").Value = 20.150343
$ws.Range("N5").Value = 40.300686
$ws.Range("O5").Value = 0.3630064529378702
$ws.Range("P5").Value = 0.2753184233934096
$ws.Range("Q5").Value = 0.5465377531889999
$ws.Range("R5").Value = 3.279226519134
$ws.Range("S5").Value = 0.03487864883747752
$ws.Range("T5").Value = 0.02645334409432725
$ws.Range("G6").Value = 0.027123
$ws.Range("H6").Value = 0.081369
$ws.Range("I6").Value = 0.0960827240265261
$ws.Range("J6").Value = 0.09608272402652611
$ws.Range("M6").Value = 0.1376283333333333
$ws.Range("N6").Value = 0.412885
$ws.Range("O6").Value = 0.002479360927359111
$ws.Range("P6").Value = 0.002820667798130978
$ws.Range("Q6").Value = 0.003732893284999999
$ws.Range("R6").Value = 0.033596039565
$ws.Range("S6").Value = 0.0002382237517455973
$ws.Range("T6").Value = 0.0002710174456183279
$ws.Range("I7").Value = 0.9039172759734738
$ws.Range("J7").Value = 0.9039172759734738
$ws.Range("M7").Value = 7.60267
$ws.Range("N7").Value = 22.80801
$ws.Range("O7").Value = 0.1369613544323865
$ws.Range("P7").Value = 0.1558153465164618
$ws.Range("Q7").Value = 1.93993529055
$ws.Range("R7").Value = 17.45941761495
$ws.Range("S7").Value = 0.1238017344121602
$ws.Range("T7").Value = 0.1408441835780231
$ws.Range("I8").Value = 0.9039172759734738
$ws.Range("J8").Value = 0.9039172759734738
$ws.Range("N8").Value = 81.75477699999999
$ws.Range("O8").Value = 0.4909347632361489
$ws.Range("P8").Value = 0.5585164557377457
$ws.Range("R8").Value = 62.58287301961499
$ws.Range("S8").Value = 0.4437644138651021
$ws.Range("T8").Value = 0.5048526732568224
$ws.Range("I9").Value = 0.9039172759734738
$ws.Range("J9").Value = 0.9039172759734738
$ws.Range("M9").Value = 0.3673663333333333
$ws.Range("N9").Value = 1.102099
$ws.Range("O9").Value = 0.00661806846623527
$ws.Range("P9").Value = 0.007529106554252038
$ws.Range("Q9").Value = 0.09373903044499998
$ws.Range("R9").Value = 0.8436512740049998
$ws.Range("S9").Value = 0.005982186420205331
$ws.Range("T9").Value = 0.00680568948703353
$ws.Range("I10").Value = 0.9039172759734738
$ws.Range("J10").Value = 0.9039172759734738
$ws.Range("M10").Value = 20.150343
$ws.Range("N10").Value = 40.300686
$ws.Range("O10").Value = 0.3630064529378702
$ws.Range("P10").Value = 0.2753184233934096
$ws.Range("Q10").Value = 5.141662271594999
$ws.Range("R10").Value = 30.84997362956999
$ws.Range("S10").Value = 0.3281278041003927
$ws.Range("T10").Value = 0.2488650792990824
$ws.Range("I11").Value = 0.9039172759734738
$ws.Range("J11").Value = 0.9039172759734738
$ws.Range("M11").Value = 0.1376283333333333
$ws.Range("N11").Value = 0.412885
$ws.Range("O11").Value = 0.002479360927359111
$ws.Range("P11").Value = 0.002820667798130978
$ws.Range("Q11").Value = 0.035117933675
$ws.Range("R11").Value = 0.316061403075
$ws.Range("S11").Value = 0.002241137175613514
$ws.Range("T11").Value = 0.00254965035251265
